$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.472.64"
$ws.Range("E2").Value = "  -0.39%  "

$ws.Range("D3").Value = "3.431.39"
$ws.Range("E3").Value = "  -2.11%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.55"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.47"
$ws.Range("E6").Value = "  -4.40%  "

$ws.Range("D7").Value = "3.430.66"
$ws.Range("E7").Value = "  -2.01%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("E10").Value = "  -4.74%  "

$ws.Range("E11").Value = "  -8.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -6.45%  "

$ws.Range("D13").Value = "4.017.56"
$ws.Range("E13").Value = "  -2.13%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("E14").Value = "  -9.41%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.40"
$ws.Range("E15").Value = "  -7.74%  "

$ws.Range("D16").Value = "3.434.57"
$ws.Range("E16").Value = "  -2.31%  "

$ws.Range("D17").Value = "65.462.77"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.86"
$ws.Range("E19").Value = "  -9.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.87"
$ws.Range("E20").Value = "  -4.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("E21").Value = "  -4.66%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.39"
$ws.Range("E22").Value = "  -4.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.555"
$ws.Range("E23").Value = "  -6.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.33"
$ws.Range("E24").Value = "  -5.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("D26").Value = "3.571.62"
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000107"
$ws.Range("E27").Value = "  -7.11%  "

$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("E29").Value = "  -5.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("E30").Value = "  -8.81%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.23"
$ws.Range("E31").Value = "  -8.33%  "

$ws.Range("D32").Value = "3.440.05"
$ws.Range("E32").Value = "  -1.93%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.145"
$ws.Range("E34").Value = "  -4.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.00"
$ws.Range("E35").Value = "  -4.69%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.45"
$ws.Range("E36").Value = "  -0.98%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.88"
$ws.Range("E37").Value = "  -7.61%  "

$ws.Range("E38").Value = "  -5.68%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.47"
$ws.Range("E39").Value = "  -5.69%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.79"
$ws.Range("E40").Value = "  -7.80%  "

$ws.Range("E41").Value = "  -6.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.822"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.66"
$ws.Range("E43").Value = "  -3.28%  "

$ws.Range("E44").Value = "  -0.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("E45").Value = "  -12.08%  "

$ws.Range("E46").Value = "  -7.98%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.12"
$ws.Range("E47").Value = "  +3.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.60"
$ws.Range("E48").Value = "  +0.88%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.53"
$ws.Range("E49").Value = "  -7.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.10"
$ws.Range("E50").Value = "  -12.06%  "

$ws.Range("D51").Value = "2.206.36"
$ws.Range("E51").Value = "  -6.16%  "
